# Realestate Update resale numbers 2023-06-28 22:36
# Appends a new data row (row 88) to the CityResaleNum sheet with the
# latest resale-number snapshot, mirroring the format of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 88

# Columns A, B and D hold date-looking / time-looking / numeric-looking
# text (e.g. "2023-06-28", "22:35:53", "26"). Excel's COM layer would
# normally auto-coerce such strings into date serials / numbers, so force
# a text number format immediately before assigning the value, then clear
# the formatting again afterwards so the cell is left, like all the other
# data rows, without any explicit style while still holding the literal
# text value.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2023-06-28"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").NumberFormat = "@"
$ws.Range("B$row").Value = "22:35:53"
$ws.Range("B$row").ClearFormats()

$ws.Range("C$row").Value = "Wednesday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "26"
$ws.Range("D$row").ClearFormats()

# Columns E-T hold the numeric resale counts per city.
$ws.Range("E$row").Value = 123107
$ws.Range("F$row").Value = 134332
$ws.Range("G$row").Value = 163951
$ws.Range("H$row").Value = 134121
$ws.Range("I$row").Value = 177220
$ws.Range("J$row").Value = 115010
$ws.Range("K$row").Value = 204391
$ws.Range("L$row").Value = 226534
$ws.Range("M$row").Value = 176294
$ws.Range("N$row").Value = 104482
$ws.Range("O$row").Value = 39773
$ws.Range("P$row").Value = 33731
$ws.Range("Q$row").Value = 52461
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36142
$ws.Range("T$row").Value = -1
